# Update pinmap with pullup todo
# - New ANALOG_LEFT / ANALOG_RIGHT rows for IO32 / IO33
# - SD_CS / IO_CS / IO_IRQ usages shift down to IO25 / IO26 / IO27
# - Several "(use pullup?)" / "(add pullup?)" annotations on CS lines
# - "IO (MCP23S17)" marked done, new "ANALOG done" section row added

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E9").Value  = "ANALOG_LEFT (ADC1_CH4 (GPIO 32))"
$ws.Range("E10").Value = "ANALOG_RIGHT (ADC1_CH5 (GPIO 33))"
$ws.Range("E11").Value = "SD_CS (use pullup?)"
$ws.Range("E12").Value = "IO_CS (use pullup?)"
$ws.Range("E13").Value = "IO_IRQ"

$ws.Range("E24").Value = "LCD_PIN_CS (use internal pullup after boot?)"
$ws.Range("E26").Value = "TOUCH_PIN_CS (use internal pullup after boot?)"
$ws.Range("E28").Value = "PSRAM_CS (add pullup?)"

$ws.Range("D48").Value = "IO (MCP23S17) done"
$ws.Range("D50").Value = "ANALOG done"

# E10 gets a wrap-text style (matches the longer combined label)
$ws.Range("E10").ClearFormats()
$ws.Range("E10").WrapText = $true

# Move the active selection like the author left it
$ws.Range("E35").Select()
